$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "39.952.47"
$ws.Range("E2").Value = "  -0.19%  "

Set-TextValue "D3" "2.204.62"
$ws.Range("E3").Value = "  -1.30%  "

$ws.Range("E4").Value = "  -0.03%  "

Set-TextValue "D5" "293.84"
$ws.Range("E5").Value = "  +0.08%  "

Set-TextValue "D6" "87.28"
$ws.Range("E6").Value = "  +0.74%  "

Set-TextValue "D7" "0.508"
$ws.Range("E7").Value = "  -1.30%  "

$ws.Range("E8").Value = "  -0.06%  "

Set-TextValue "D9" "0.467"
$ws.Range("E9").Value = "  -0.82%  "

Set-TextValue "D10" "29.96"
$ws.Range("E10").Value = "  -2.76%  "

Set-TextValue "D11" "0.0776"
$ws.Range("E11").Value = "  -1.75%  "

Set-TextValue "D12" "49.90"
$ws.Range("E12").Value = "  +5.96%  "

Set-TextValue "D13" "0.111"
$ws.Range("E13").Value = "  +2.56%  "

Set-TextValue "D14" "6.42"
$ws.Range("E14").Value = "  +0.23%  "

Set-TextValue "D15" "2.547.71"
$ws.Range("E15").Value = "  -1.12%  "

Set-TextValue "D16" "13.68"
$ws.Range("E16").Value = "  -2.89%  "

Set-TextValue "D17" "2.196.51"
$ws.Range("E17").Value = "  -1.07%  "

Set-TextValue "D18" "0.723"
$ws.Range("E18").Value = "  -0.99%  "

Set-TextValue "D19" "39.853.03"
$ws.Range("E19").Value = "  -0.24%  "

Set-TextValue "D20" "0.0₃0882"
$ws.Range("E20").Value = "  -0.78%  "

Set-TextValue "D21" "11.21"
$ws.Range("E21").Value = "  +2.02%  "

Set-TextValue "D22" "5.75"
$ws.Range("E22").Value = "  -0.90%  "

Set-TextValue "D23" "65.11"
$ws.Range("E23").Value = "  -0.24%  "

Set-TextValue "D24" "236.51"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("E25").Value = "  +0.05%  "

$ws.Range("E26").Value = "  -0.54%  "

Set-TextValue "D27" "1.80"
$ws.Range("E27").Value = "  -2.60%  "

Set-TextValue "D28" "22.48"
$ws.Range("E28").Value = "  -1.41%  "

Set-TextValue "D29" "2.06"
$ws.Range("E29").Value = "  -7.39%  "

Set-TextValue "D30" "9.14"
$ws.Range("E30").Value = "  -1.51%  "

Set-TextValue "D31" "156.67"
$ws.Range("E31").Value = "  +2.57%  "

Set-TextValue "D32" "31.22"
$ws.Range("E32").Value = "  -5.80%  "

$ws.Range("E33").Value = "  -0.02%  "

Set-TextValue "D34" "4.87"
$ws.Range("E34").Value = "  -1.19%  "

Set-TextValue "D35" "0.0708"
$ws.Range("E35").Value = "  -1.74%  "

$ws.Range("E36").Value = "  -2.37%  "

Set-TextValue "D37" "2.81"
$ws.Range("E37").Value = "  +0.80%  "

$ws.Range("E38").Value = "  +0.79%  "

Set-TextValue "D39" "0.0973"
$ws.Range("E39").Value = "  -2.60%  "

Set-TextValue "D40" "15.18"
$ws.Range("E40").Value = "  -6.27%  "

$ws.Range("E41").Value = "  -1.98%  "

Set-TextValue "D42" "2.115.67"
$ws.Range("E42").Value = "  +3.71%  "

Set-TextValue "D43" "3.71"
$ws.Range("E43").Value = "  -2.64%  "

Set-TextValue "D44" "2.12"
$ws.Range("E44").Value = "  -2.06%  "

Set-TextValue "D45" "0.0266"
$ws.Range("E45").Value = "  -1.41%  "

Set-TextValue "D46" "17.31"
$ws.Range("E46").Value = "  +3.98%  "

Set-TextValue "D47" "9.65"
$ws.Range("E47").Value = "  -4.14%  "

Set-TextValue "D48" "2.64"
$ws.Range("E48").Value = "  +2.85%  "

Set-TextValue "D49" "2.421.03"
$ws.Range("E49").Value = "  -1.67%  "

Set-TextValue "D50" "1.47"
$ws.Range("E50").Value = "  +2.81%  "

$ws.Range("E51").Value = "  +1.42%  "
